# Auto-generated edit script applying the Ultima_Profits.xlsx diff
# Updates currentAveragePrice / currentAveragePriceNQ/HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ
# columns (H..N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15924
$ws.Range("H54").Value = 3489.5
$ws.Range("J54").Value = 4980
$ws.Range("L54").Value = 4980
$ws.Range("N54").Value = -5952
$ws.Range("H76").Value = 3238.8948
$ws.Range("I76").Value = 3163.3333
$ws.Range("J76").Value = 3522.25
$ws.Range("K76").Value = 3163.3333
$ws.Range("L76").Value = 3522.25
$ws.Range("M76").Value = -2848.3333
$ws.Range("N76").Value = -4152.25
$ws.Range("H79").Value = 3238.8948
$ws.Range("I79").Value = 3163.3333
$ws.Range("J79").Value = 3522.25
$ws.Range("K79").Value = 3163.3333
$ws.Range("L79").Value = 3522.25
$ws.Range("M79").Value = -2071.3333
$ws.Range("N79").Value = -5706.25
$ws.Range("H100").Value = 2241.1428
$ws.Range("J100").Value = 3203.6
$ws.Range("L100").Value = 3203.6
$ws.Range("N100").Value = -4285.6
$ws.Range("H101").Value = 1501.5
$ws.Range("I101").Value = 1501.5
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 4504.5
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -2882.5
$ws.Range("N101").Value = ""
$ws.Range("H127").Value = 1385.2916
$ws.Range("I127").Value = 394.7
$ws.Range("K127").Value = 1184.1
$ws.Range("M127").Value = 3775.9
$ws.Range("H129").Value = 1252.4889
$ws.Range("I129").Value = 354.93332
$ws.Range("J129").Value = 1432
$ws.Range("K129").Value = 1064.79996
$ws.Range("L129").Value = 4296
$ws.Range("M129").Value = 3935.20004
$ws.Range("N129").Value = -14296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6048.48
$ws.Range("I32").Value = 3308.7036
$ws.Range("J32").Value = 17728.578
$ws.Range("K32").Value = 3308.7036
$ws.Range("L32").Value = 17728.578
$ws.Range("M32").Value = -3021.7036
$ws.Range("N32").Value = -18302.578
$ws.Range("H61").Value = 10640873
$ws.Range("I61").Value = 13516230
$ws.Range("K61").Value = 13516230
$ws.Range("M61").Value = -13516018
$ws.Range("H74").Value = 27780700
$ws.Range("I74").Value = 83335550
$ws.Range("J74").Value = 3271.1667
$ws.Range("K74").Value = 83335550
$ws.Range("L74").Value = 3271.1667
$ws.Range("M74").Value = -83334676
$ws.Range("N74").Value = -5019.1667
$ws.Range("H75").Value = 44086.5
$ws.Range("J75").Value = 44086.5
$ws.Range("L75").Value = 44086.5
$ws.Range("N75").Value = -45834.5
$ws.Range("H77").Value = 27780700
$ws.Range("I77").Value = 83335550
$ws.Range("J77").Value = 3271.1667
$ws.Range("K77").Value = 416677750
$ws.Range("L77").Value = 16355.8335
$ws.Range("M77").Value = -416673382
$ws.Range("N77").Value = -25091.8335
$ws.Range("H78").Value = 44086.5
$ws.Range("J78").Value = 44086.5
$ws.Range("L78").Value = 132259.5
$ws.Range("N78").Value = -140995.5
$ws.Range("H132").Value = 7144647.5
$ws.Range("I132").Value = 10418259
$ws.Range("K132").Value = 31254777
$ws.Range("M132").Value = -31252247
$ws.Range("H136").Value = 10640873
$ws.Range("I136").Value = 13516230
$ws.Range("K136").Value = 40548690
$ws.Range("M136").Value = -40546140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 33143
$ws.Range("I26").Value = 20111
$ws.Range("K26").Value = 20111
$ws.Range("M26").Value = -19819
$ws.Range("H134").Value = 5526.8335
$ws.Range("I134").Value = 4911.6
$ws.Range("J134").Value = 5966.2856
$ws.Range("K134").Value = 14734.8
$ws.Range("L134").Value = 17898.8568
$ws.Range("M134").Value = -12199.8
$ws.Range("N134").Value = -22968.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4764911.5
$ws.Range("I31").Value = 2891.3115
$ws.Range("J31").Value = 37040828
$ws.Range("K31").Value = 2891.3115
$ws.Range("L31").Value = 37040828
$ws.Range("M31").Value = -2596.3115
$ws.Range("N31").Value = -37041418
$ws.Range("H34").Value = 4764911.5
$ws.Range("I34").Value = 2891.3115
$ws.Range("J34").Value = 37040828
$ws.Range("K34").Value = 2891.3115
$ws.Range("L34").Value = 37040828
$ws.Range("M34").Value = -2689.3115
$ws.Range("N34").Value = -37041232
$ws.Range("H58").Value = 2030.6364
$ws.Range("I58").Value = 1298.5883
$ws.Range("J58").Value = 2808.4375
$ws.Range("K58").Value = 1298.5883
$ws.Range("L58").Value = 2808.4375
$ws.Range("M58").Value = -1095.5883
$ws.Range("N58").Value = -3214.4375
$ws.Range("H88").Value = 44383.285
$ws.Range("J88").Value = 48395.332
$ws.Range("L88").Value = 48395.332
$ws.Range("N88").Value = -49207.332
$ws.Range("H91").Value = 44383.285
$ws.Range("J91").Value = 48395.332
$ws.Range("L91").Value = 48395.332
$ws.Range("N91").Value = -51203.332
$ws.Range("H99").Value = 1531.6538
$ws.Range("I99").Value = 1366.4286
$ws.Range("J99").Value = 2225.6
$ws.Range("K99").Value = 1366.4286
$ws.Range("L99").Value = 2225.6
$ws.Range("M99").Value = 131.5714
$ws.Range("N99").Value = -5221.6
$ws.Range("H106").Value = 47670.75
$ws.Range("J106").Value = 47670.75
$ws.Range("L106").Value = 47670.75
$ws.Range("N106").Value = -50194.75
$ws.Range("H126").Value = 1531.6538
$ws.Range("I126").Value = 1366.4286
$ws.Range("J126").Value = 2225.6
$ws.Range("K126").Value = 4099.2858
$ws.Range("L126").Value = 6676.799999999999
$ws.Range("M126").Value = -1629.2858
$ws.Range("N126").Value = -11616.8
$ws.Range("H132").Value = 31253928
$ws.Range("I132").Value = 71432720
$ws.Range("J132").Value = 3758.3333
$ws.Range("K132").Value = 214298160
$ws.Range("L132").Value = 11274.9999
$ws.Range("M132").Value = -214295630
$ws.Range("N132").Value = -16334.9999
$ws.Range("H134").Value = 3407935.5
$ws.Range("I134").Value = 11341.333
$ws.Range("J134").Value = 5955381
$ws.Range("K134").Value = 34023.999
$ws.Range("L134").Value = 17866143
$ws.Range("M134").Value = -31488.999
$ws.Range("N134").Value = -17871213
$ws.Range("H136").Value = 2030.6364
$ws.Range("I136").Value = 1298.5883
$ws.Range("J136").Value = 2808.4375
$ws.Range("K136").Value = 3895.7649
$ws.Range("L136").Value = 8425.3125
$ws.Range("M136").Value = -1345.7649
$ws.Range("N136").Value = -13525.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1495.7778
$ws.Range("I44").Value = 300
$ws.Range("J44").Value = 1645.25
$ws.Range("K44").Value = 900
$ws.Range("L44").Value = 4935.75
$ws.Range("M44").Value = -502
$ws.Range("N44").Value = -5731.75
$ws.Range("H113").Value = 50000880
$ws.Range("I113").Value = 166667220
$ws.Range("J113").Value = 1017.5714
$ws.Range("K113").Value = 500001660
$ws.Range("L113").Value = 3052.7142
$ws.Range("M113").Value = -499999490
$ws.Range("N113").Value = -7392.7142
$ws.Range("H140").Value = 4323.3125
$ws.Range("I140").Value = 2014.4166
$ws.Range("K140").Value = 6043.2498
$ws.Range("M140").Value = -863.2497999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7601.885
$ws.Range("J70").Value = 4025.1428
$ws.Range("L70").Value = 4025.1428
$ws.Range("N70").Value = -4565.1428
$ws.Range("H73").Value = 7601.885
$ws.Range("J73").Value = 4025.1428
$ws.Range("L73").Value = 4025.1428
$ws.Range("N73").Value = -5897.1428
$ws.Range("H118").Value = 17316.666
$ws.Range("J118").Value = 17316.666
$ws.Range("L118").Value = 17316.666
$ws.Range("N118").Value = -20630.666
$ws.Range("H132").Value = 4728.952
$ws.Range("I132").Value = 6115.885
$ws.Range("K132").Value = 18347.655
$ws.Range("M132").Value = -15817.655

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10422648
$ws.Range("I132").Value = 3200.3
$ws.Range("J132").Value = 27788394
$ws.Range("K132").Value = 9600.900000000001
$ws.Range("L132").Value = 83365182
$ws.Range("M132").Value = -7070.900000000001
$ws.Range("N132").Value = -83370242
$ws.Range("H136").Value = 25007470
$ws.Range("I136").Value = 27780438
$ws.Range("J136").Value = 50752.5
$ws.Range("K136").Value = 83341314
$ws.Range("L136").Value = 152257.5
$ws.Range("M136").Value = -83338764
$ws.Range("N136").Value = -157357.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2623.6924
$ws.Range("I132").Value = 2478.2222
$ws.Range("J132").Value = 2951
$ws.Range("K132").Value = 7434.6666
$ws.Range("L132").Value = 8853
$ws.Range("M132").Value = -4904.6666
$ws.Range("N132").Value = -13913
$ws.Range("H136").Value = 1105.75
$ws.Range("I136").Value = 1249.1578
$ws.Range("K136").Value = 3747.4734
$ws.Range("M136").Value = -1249.1578
